$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.081400327803546
$ws.Range("D2").Value = 1.074562766688935
$ws.Range("E2").Value = 1.095227783457438
$ws.Range("F2").Value = 1.102538653272681
$ws.Range("I2").Value = 1.062478403897568
$ws.Range("J2").Value = 1.08627509629059
$ws.Range("K2").Value = 1.077251649604056
$ws.Range("L2").Value = 1.097863121809419
$ws.Range("M2").Value = 1.105155566555029
$ws.Range("N2").Value = 1.087817731524818
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.083042002065808
$ws.Range("D3").Value = 1.07585574606842
$ws.Range("E3").Value = 1.096900055763648
$ws.Range("F3").Value = 1.104353676254857
$ws.Range("I3").Value = 1.063118674605324
$ws.Range("J3").Value = 1.087575281359329
$ws.Range("K3").Value = 1.078360610103977
$ws.Range("L3").Value = 1.099354377496846
$ws.Range("M3").Value = 1.106790592075053
$ws.Range("N3").Value = 1.089119763005488
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.084100554421431
$ws.Range("D4").Value = 1.076688553059084
$ws.Range("E4").Value = 1.097978731247905
$ws.Range("F4").Value = 1.105524910330211
$ws.Range("I4").Value = 1.063529521904468
$ws.Range("J4").Value = 1.088412467652522
$ws.Range("K4").Value = 1.079073795846287
$ws.Range("L4").Value = 1.10031544142448
$ws.Range("M4").Value = 1.10784489827166
$ws.Range("N4").Value = 1.089958138199244
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.084544694870168
$ws.Range("D5").Value = 1.07703775909934
$ws.Range("E5").Value = 1.098431408788879
$ws.Range("F5").Value = 1.106016545979018
$ws.Range("I5").Value = 1.06370142346644
$ws.Range("J5").Value = 1.088763447575424
$ws.Range("K5").Value = 1.079372581127689
$ws.Range("L5").Value = 1.100718558796105
$ws.Range("M5").Value = 1.108287266845303
$ws.Range("N5").Value = 1.090309616553894
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.084619217153021
$ws.Range("D6").Value = 1.077096339564283
$ws.Range("E6").Value = 1.098507369013822
$ws.Range("F6").Value = 1.106099050200132
$ws.Range("I6").Value = 1.063730238704541
$ws.Range("J6").Value = 1.088822321962225
$ws.Range("K6").Value = 1.07942268796775
$ws.Range("L6").Value = 1.100786190804854
$ws.Range("M6").Value = 1.10836149239136
$ws.Range("N6").Value = 1.090368574549078
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.084106492464315
$ws.Range("D7").Value = 1.076693222710094
$ws.Range("E7").Value = 1.097984783060344
$ws.Range("F7").Value = 1.105531482518624
$ws.Range("I7").Value = 1.063531822068063
$ws.Range("J7").Value = 1.088417161269343
$ws.Range("K7").Value = 1.079077792292125
$ws.Range("L7").Value = 1.100320831469675
$ws.Range("M7").Value = 1.107850812587025
$ws.Range("N7").Value = 1.089962838481539
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.081955919236475
$ws.Range("D8").Value = 1.075000537247637
$ws.Range("E8").Value = 1.095793648360611
$ws.Range("F8").Value = 1.103152722918058
$ws.Range("I8").Value = 1.062695505373157
$ws.Range("J8").Value = 1.086715361081471
$ws.Range("K8").Value = 1.077627343363876
$ws.Range("L8").Value = 1.098367910463327
$ws.Range("M8").Value = 1.105708899292819
$ws.Range("N8").Value = 1.088258621542201
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.078137126520849
$ws.Range("D9").Value = 1.071987852387151
$ws.Range("E9").Value = 1.091905876030552
$ws.Range("F9").Value = 1.098935730344065
$ws.Range("I9").Value = 1.061195043760582
$ws.Range("J9").Value = 1.083684416815331
$ws.Range("K9").Value = 1.075037335621091
$ws.Range("L9").Value = 1.094896254665589
$ws.Range("M9").Value = 1.101905804632014
$ws.Range("N9").Value = 1.085223372987548
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.075570625389669
$ws.Range("D10").Value = 1.069958453502806
$ws.Range("E10").Value = 1.089295083103433
$ws.Range("F10").Value = 1.096106325653557
$ws.Range("I10").Value = 1.060176266517008
$ws.Range("J10").Value = 1.081641327041812
$ws.Range("K10").Value = 1.073286954240416
$ws.Range("L10").Value = 1.092560481165627
$ws.Range("M10").Value = 1.099350076181228
$ws.Range("N10").Value = 1.083177381792206
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.07445418214412
$ws.Range("D11").Value = 1.069074554341474
$ws.Range("E11").Value = 1.088159866788177
$ws.Range("F11").Value = 1.09487663793464
$ws.Range("I11").Value = 1.059730637433433
$ws.Range("J11").Value = 1.080751130620834
$ws.Range("K11").Value = 1.072523225287201
$ws.Range("L11").Value = 1.09154379745792
$ws.Range("M11").Value = 1.098238373351859
$ws.Range("N11").Value = 1.082285921190203
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.074038695008102
$ws.Range("D12").Value = 1.068745444571435
$ws.Range("E12").Value = 1.087737467762276
$ws.Range("F12").Value = 1.094419175026265
$ws.Range("I12").Value = 1.059564426823697
$ws.Range("J12").Value = 1.080419625559544
$ws.Range("K12").Value = 1.072238655765769
$ws.Range("L12").Value = 1.091165344847796
$ws.Range("M12").Value = 1.097824658646187
$ws.Range("N12").Value = 1.081953945353702
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.07412785449212
$ws.Range("D13").Value = 1.068816075748571
$ws.Range("E13").Value = 1.087828107092443
$ws.Range("F13").Value = 1.094517334478654
$ws.Range("I13").Value = 1.059600110672191
$ws.Range("J13").Value = 1.080490773090765
$ws.Range("K13").Value = 1.072299737254476
$ws.Range("L13").Value = 1.091246561208842
$ws.Range("M13").Value = 1.097913437407699
$ws.Range("N13").Value = 1.082025193922578
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.074419854085459
$ws.Range("D14").Value = 1.06904736623753
$ws.Range("E14").Value = 1.08812496613884
$ws.Range("F14").Value = 1.094838838351032
$ws.Range("I14").Value = 1.059716912422689
$ws.Range("J14").Value = 1.080723745672484
$ws.Range("K14").Value = 1.072499720854458
$ws.Range("L14").Value = 1.091512531096476
$ws.Range("M14").Value = 1.098204191564196
$ws.Range("N14").Value = 1.082258497352085
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.074599659348685
$ws.Range("D15").Value = 1.069189766752426
$ws.Range("E15").Value = 1.088307773599551
$ws.Range("F15").Value = 1.095036833838987
$ws.Range("I15").Value = 1.059788786902937
$ws.Range("J15").Value = 1.08086717510892
$ws.Range("K15").Value = 1.072622819422014
$ws.Range("L15").Value = 1.091676295912835
$ws.Range("M15").Value = 1.098383231020959
$ws.Range("N15").Value = 1.082402130474773
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.075644610282715
$ws.Range("D16").Value = 1.07001700501776
$ws.Range("E16").Value = 1.089370322476415
$ws.Range("F16").Value = 1.096187838645475
$ws.Range("I16").Value = 1.060205746024782
$ws.Range("J16").Value = 1.081700288645047
$ws.Range("K16").Value = 1.073337516904847
$ws.Range("L16").Value = 1.092627842374826
$ws.Range("M16").Value = 1.099423748085743
$ws.Range("N16").Value = 1.08323642712768
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.076298692920608
$ws.Range("D17").Value = 1.070534518336025
$ws.Range("E17").Value = 1.090035552654574
$ws.Range("F17").Value = 1.096908604503379
$ws.Range("I17").Value = 1.060466084365175
$ws.Range("J17").Value = 1.082221388030039
$ws.Range("K17").Value = 1.07378426389346
$ws.Range("L17").Value = 1.093223296528225
$ws.Range("M17").Value = 1.10007506993322
$ws.Range("N17").Value = 1.08375826653357
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.076679714487641
$ws.Range("D18").Value = 1.070835878680596
$ws.Range("E18").Value = 1.090423115502333
$ws.Range("F18").Value = 1.097328578696867
$ws.Range("I18").Value = 1.060617502552841
$ws.Range("J18").Value = 1.082524804514391
$ws.Range("K18").Value = 1.074044285032547
$ws.Range("L18").Value = 1.093570106802983
$ws.Range("M18").Value = 1.100454488472409
$ws.Range("N18").Value = 1.084062113904132
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.076809549829417
$ws.Range("D19").Value = 1.07093855103884
$ws.Range("E19").Value = 1.090555187875379
$ws.Range("F19").Value = 1.097471705650765
$ws.Range("I19").Value = 1.060669059116939
$ws.Range("J19").Value = 1.082628171956025
$ws.Range("K19").Value = 1.074132851182374
$ws.Range("L19").Value = 1.093688274509837
$ws.Range("M19").Value = 1.10058377845303
$ws.Range("N19").Value = 1.08416562813939
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.076228567257108
$ws.Range("D20").Value = 1.070479045514235
$ws.Range("E20").Value = 1.089964226917881
$ws.Range("F20").Value = 1.096831318428531
$ws.Range("I20").Value = 1.060438197354567
$ws.Range("J20").Value = 1.082165534140101
$ws.Range("K20").Value = 1.073736390078165
$ws.Range("L20").Value = 1.093159462613552
$ws.Range("M20").Value = 1.10000523972014
$ws.Range("N20").Value = 1.083702333324701
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.074333889459867
$ws.Range("D21").Value = 1.06897927889502
$ws.Range("E21").Value = 1.08803756883157
$ws.Range("F21").Value = 1.094744183022524
$ws.Range("I21").Value = 1.059682536178766
$ws.Range("J21").Value = 1.080655164553427
$ws.Range("K21").Value = 1.072440855263045
$ws.Range("L21").Value = 1.091434232127322
$ws.Range("M21").Value = 1.098118593365725
$ws.Range("N21").Value = 1.082189818839973
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.073138047458521
$ws.Range("D22").Value = 1.06803173597237
$ws.Range("E22").Value = 1.086821974175306
$ws.Range("F22").Value = 1.093427845442131
$ws.Range("I22").Value = 1.059203459174151
$ws.Range("J22").Value = 1.079700629531715
$ws.Range("K22").Value = 1.071621164108556
$ws.Range("L22").Value = 1.090344809890701
$ws.Range("M22").Value = 1.096927868454614
$ws.Range("N22").Value = 1.081233928269044
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.073772426901671
$ws.Range("D23").Value = 1.068534485816731
$ws.Range("E23").Value = 1.08746679114873
$ws.Range("F23").Value = 1.094126053868194
$ws.Range("I23").Value = 1.059457805719982
$ws.Range("J23").Value = 1.080207117141256
$ws.Range("K23").Value = 1.072056189837144
$ws.Range("L23").Value = 1.090922784933968
$ws.Range("M23").Value = 1.097559528759963
$ws.Range("N23").Value = 1.081741135149085
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.076260255551258
$ws.Range("D24").Value = 1.070504112828224
$ws.Range("E24").Value = 1.089996457352634
$ws.Range("F24").Value = 1.096866242026887
$ws.Range("I24").Value = 1.060450799631297
$ws.Range("J24").Value = 1.082190773752195
$ws.Range("K24").Value = 1.073758023921011
$ws.Range("L24").Value = 1.093188307981418
$ws.Range("M24").Value = 1.100036794490435
$ws.Range("N24").Value = 1.08372760877994
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.07912793962479
$ws.Range("D25").Value = 1.072770335522651
$ws.Range("E25").Value = 1.092914224603371
$ws.Range("F25").Value = 1.100029031851323
$ws.Range("I25").Value = 1.061586169278816
$ws.Range("J25").Value = 1.084471885031795
$ws.Range("K25").Value = 1.075711037265888
$ws.Range("L25").Value = 1.095797453116983
$ws.Range("M25").Value = 1.102892506188762
$ws.Range("N25").Value = 1.086011959499204
